# Auto-generated edit script applying the Gungnir_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N)
# for specific rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 1601.8148
$ws.Cells.Item(15, 9).Value = 1601.8148
$ws.Cells.Item(15, 11).Value = 4805.4444
$ws.Cells.Item(15, 13).Value = -4636.4444
# Row 19
$ws.Cells.Item(19, 8).Value = 774.7143
$ws.Cells.Item(19, 9).Value = 940
$ws.Cells.Item(19, 10).Value = 723.0625
$ws.Cells.Item(19, 11).Value = 940
$ws.Cells.Item(19, 12).Value = 723.0625
$ws.Cells.Item(19, 13).Value = -765
$ws.Cells.Item(19, 14).Value = -1073.0625
# Row 54
$ws.Cells.Item(54, 8).Value = 6500
$ws.Cells.Item(54, 9).Value = 5000
$ws.Cells.Item(54, 11).Value = 5000
$ws.Cells.Item(54, 13).Value = -4514
# Row 138
$ws.Cells.Item(138, 8).Value = 2309.4412
$ws.Cells.Item(138, 9).Value = 1532
$ws.Cells.Item(138, 10).Value = 3294.2
$ws.Cells.Item(138, 11).Value = 4596
$ws.Cells.Item(138, 12).Value = 9882.599999999999
$ws.Cells.Item(138, 13).Value = 544
$ws.Cells.Item(138, 14).Value = -20162.6

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 14088374
$ws.Cells.Item(32, 9).Value = 3844.3833
$ws.Cells.Item(32, 11).Value = 3844.3833
$ws.Cells.Item(32, 13).Value = -3557.3833
# Row 61
$ws.Cells.Item(61, 8).Value = 4066073.2
$ws.Cells.Item(61, 9).Value = 4505553
$ws.Cells.Item(61, 10).Value = 887.5
$ws.Cells.Item(61, 11).Value = 4505553
$ws.Cells.Item(61, 12).Value = 887.5
$ws.Cells.Item(61, 13).Value = -4505341
$ws.Cells.Item(61, 14).Value = -1311.5
# Row 68
$ws.Cells.Item(68, 8).Value = 19776
$ws.Cells.Item(68, 10).Value = 19776
$ws.Cells.Item(68, 12).Value = 19776
$ws.Cells.Item(68, 14).Value = -21398
# Row 71
$ws.Cells.Item(71, 8).Value = 19776
$ws.Cells.Item(71, 10).Value = 19776
$ws.Cells.Item(71, 12).Value = 59328
$ws.Cells.Item(71, 14).Value = -67440
# Row 75
$ws.Cells.Item(75, 8).Value = 25386.5
$ws.Cells.Item(75, 10).Value = 25386.5
$ws.Cells.Item(75, 12).Value = 25386.5
$ws.Cells.Item(75, 14).Value = -27134.5
# Row 78
$ws.Cells.Item(78, 8).Value = 25386.5
$ws.Cells.Item(78, 10).Value = 25386.5
$ws.Cells.Item(78, 12).Value = 76159.5
$ws.Cells.Item(78, 14).Value = -84895.5
# Row 94
$ws.Cells.Item(94, 8).Value = 43456
$ws.Cells.Item(94, 10).Value = 43456
$ws.Cells.Item(94, 12).Value = 43456
$ws.Cells.Item(94, 14).Value = -45258
# Row 132
$ws.Cells.Item(132, 8).Value = 1051296.6
$ws.Cells.Item(132, 9).Value = 706.29785
$ws.Cells.Item(132, 10).Value = 6537713
$ws.Cells.Item(132, 11).Value = 2118.89355
$ws.Cells.Item(132, 12).Value = 19613139
$ws.Cells.Item(132, 13).Value = 411.1064499999998
$ws.Cells.Item(132, 14).Value = -19618199
# Row 136
$ws.Cells.Item(136, 8).Value = 4066073.2
$ws.Cells.Item(136, 9).Value = 4505553
$ws.Cells.Item(136, 10).Value = 887.5
$ws.Cells.Item(136, 11).Value = 13516659
$ws.Cells.Item(136, 12).Value = 2662.5
$ws.Cells.Item(136, 13).Value = -13514109
$ws.Cells.Item(136, 14).Value = -7762.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 2635.6667
$ws.Cells.Item(94, 9).Value = 2972.5
$ws.Cells.Item(94, 10).Value = 2366.2
$ws.Cells.Item(94, 11).Value = 2972.5
$ws.Cells.Item(94, 12).Value = 2366.2
$ws.Cells.Item(94, 13).Value = -2521.5
$ws.Cells.Item(94, 14).Value = -3268.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 80
$ws.Cells.Item(80, 8).Value = 17000
$ws.Cells.Item(80, 10).Value = 17000
$ws.Cells.Item(80, 12).Value = 17000
$ws.Cells.Item(80, 14).Value = -19246
# Row 81
$ws.Cells.Item(81, 8).Value = 52220
$ws.Cells.Item(81, 10).Value = 52220
$ws.Cells.Item(81, 12).Value = 52220
$ws.Cells.Item(81, 14).Value = -54216
# Row 83
$ws.Cells.Item(83, 8).Value = 17000
$ws.Cells.Item(83, 10).Value = 17000
$ws.Cells.Item(83, 12).Value = 51000
$ws.Cells.Item(83, 14).Value = -62232
# Row 84
$ws.Cells.Item(84, 8).Value = 52220
$ws.Cells.Item(84, 10).Value = 52220
$ws.Cells.Item(84, 12).Value = 156660
$ws.Cells.Item(84, 14).Value = -166644
# Row 134
$ws.Cells.Item(134, 8).Value = 21740184
$ws.Cells.Item(134, 9).Value = 1049.1177
$ws.Cells.Item(134, 10).Value = 83334400
$ws.Cells.Item(134, 11).Value = 3147.3531
$ws.Cells.Item(134, 12).Value = 250003200
$ws.Cells.Item(134, 13).Value = -612.3531000000003
$ws.Cells.Item(134, 14).Value = -250008270

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 112
$ws.Cells.Item(112, 8).Value = 111113970
$ws.Cells.Item(112, 10).Value = 166670450
$ws.Cells.Item(112, 12).Value = 500011350
$ws.Cells.Item(112, 14).Value = -500013566
# Row 137
$ws.Cells.Item(137, 8).Value = 12501299
$ws.Cells.Item(137, 9).Value = 27778468
$ws.Cells.Item(137, 10).Value = 1797.7273
$ws.Cells.Item(137, 11).Value = 83335404
$ws.Cells.Item(137, 12).Value = 5393.1819
$ws.Cells.Item(137, 13).Value = -83330304
$ws.Cells.Item(137, 14).Value = -15593.1819

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Cells.Item(64, 8).Value = 14999.667
$ws.Cells.Item(64, 10).Value = 14999.667
$ws.Cells.Item(64, 12).Value = 14999.667
$ws.Cells.Item(64, 14).Value = -15495.667
# Row 67
$ws.Cells.Item(67, 8).Value = 14999.667
$ws.Cells.Item(67, 10).Value = 14999.667
$ws.Cells.Item(67, 12).Value = 14999.667
$ws.Cells.Item(67, 14).Value = -16715.667
# Row 70
$ws.Cells.Item(70, 8).Value = 11872
$ws.Cells.Item(70, 9).Value = 12400.615
$ws.Cells.Item(70, 10).Value = 5000
$ws.Cells.Item(70, 11).Value = 12400.615
$ws.Cells.Item(70, 12).Value = 5000
$ws.Cells.Item(70, 13).Value = -12130.615
$ws.Cells.Item(70, 14).Value = -5540
# Row 73
$ws.Cells.Item(73, 8).Value = 11872
$ws.Cells.Item(73, 9).Value = 12400.615
$ws.Cells.Item(73, 10).Value = 5000
$ws.Cells.Item(73, 11).Value = 12400.615
$ws.Cells.Item(73, 12).Value = 5000
$ws.Cells.Item(73, 13).Value = -11464.615
$ws.Cells.Item(73, 14).Value = -6872
# Row 95
$ws.Cells.Item(95, 8).Value = 30500
$ws.Cells.Item(95, 10).Value = 30500
$ws.Cells.Item(95, 12).Value = 30500
$ws.Cells.Item(95, 14).Value = -35992
# Row 122
$ws.Cells.Item(122, 8).Value = 55567760
$ws.Cells.Item(122, 9).Value = 62513604
$ws.Cells.Item(122, 10).Value = 1000
$ws.Cells.Item(122, 11).Value = 187540812
$ws.Cells.Item(122, 12).Value = 3000
$ws.Cells.Item(122, 13).Value = -187538362
$ws.Cells.Item(122, 14).Value = -7900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 2698
$ws.Cells.Item(22, 9).Value = 10000
$ws.Cells.Item(22, 10).Value = 1785.25
$ws.Cells.Item(22, 11).Value = 10000
$ws.Cells.Item(22, 12).Value = 1785.25
$ws.Cells.Item(22, 13).Value = -9705
$ws.Cells.Item(22, 14).Value = -2375.25
# Row 27
$ws.Cells.Item(27, 8).Value = 2698
$ws.Cells.Item(27, 9).Value = 10000
$ws.Cells.Item(27, 10).Value = 1785.25
$ws.Cells.Item(27, 11).Value = 10000
$ws.Cells.Item(27, 12).Value = 1785.25
$ws.Cells.Item(27, 13).Value = -9893
$ws.Cells.Item(27, 14).Value = -1999.25
# Row 62
$ws.Cells.Item(62, 8).Value = 7139
$ws.Cells.Item(62, 10).Value = 5122
$ws.Cells.Item(62, 12).Value = 5122
$ws.Cells.Item(62, 14).Value = -6370
# Row 65
$ws.Cells.Item(65, 8).Value = 7139
$ws.Cells.Item(65, 10).Value = 5122
$ws.Cells.Item(65, 12).Value = 15366
$ws.Cells.Item(65, 14).Value = -21606
